# Update odds values on Sheet1 as per the 2024-10-16 FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 (Danubio - Penarol)
$ws.Range("M11").Value = 1.08
$ws.Range("O11").Value = 1.44
$ws.Range("P11").Value = 2.63

# Row 12 (Boston River - Progreso)
$ws.Range("M12").Value = 1.07
$ws.Range("O12").Value = 1.33

# Row 13 (River Plate - Defensor Sp.)
$ws.Range("M13").Value = 1.06
$ws.Range("O13").Value = 1.3
$ws.Range("Q13").Value = 2.05
$ws.Range("R13").Value = 1.8
